$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 mirrors the layout of row 5 (a new statistics record appended
# below the existing A1:DF5 data block). Numeric cells first:
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 75600
$ws.Range("C6").Value = 416
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 1
$ws.Range("AC6").Value = 0

# K6:AB6 hold literal empty-string text cells (same as K5:AB5). A bare
# Value/Formula assignment of "" clears a cell instead of leaving an empty
# text value behind, so seed them as quote-prefixed empty text ("'") and
# then strip the resulting quote-prefix formatting, leaving a plain empty
# text cell.
$ws.Range("K6:AB6").Formula = "'"
$ws.Range("K6:AB6").ClearFormats()
